# The authored change repoints the deck's theme: the "Integral" color
# palette that currently drives the slide master (ppt/theme/theme2.xml,
# the theme actually wired to Presentation.SlideMaster) is replaced by
# the default "Office Theme" palette (the palette that used to live,
# unused, in ppt/theme/theme1.xml), while the font scheme / format
# scheme (already identical between the two theme parts) are untouched.
#
# We drive this through the real Theme object model exposed on the
# slide master: SlideMaster.Theme.ThemeColorScheme is the 12-slot
# scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) backing
# <a:clrScheme> in the theme part, and each entry's .RGB is settable,
# same as PowerPoint's real COM automation surface.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Target values below are the standard Office default theme colors
# (the palette theme1.xml already carried) expressed as COM RGB()
# integers: r + g*256 + b*65536.

# 1  dk1      000000
$colors.Item(1).RGB  = 0x00 + (0x00 * 256) + (0x00 * 65536)
# 2  lt1      FFFFFF
$colors.Item(2).RGB  = 0xFF + (0xFF * 256) + (0xFF * 65536)
# 3  dk2      44546A
$colors.Item(3).RGB  = 0x44 + (0x54 * 256) + (0x6A * 65536)
# 4  lt2      E7E6E6
$colors.Item(4).RGB  = 0xE7 + (0xE6 * 256) + (0xE6 * 65536)
# 5  accent1  5B9BD5
$colors.Item(5).RGB  = 0x5B + (0x9B * 256) + (0xD5 * 65536)
# 6  accent2  ED7D31
$colors.Item(6).RGB  = 0xED + (0x7D * 256) + (0x31 * 65536)
# 7  accent3  A5A5A5
$colors.Item(7).RGB  = 0xA5 + (0xA5 * 256) + (0xA5 * 65536)
# 8  accent4  FFC000
$colors.Item(8).RGB  = 0xFF + (0xC0 * 256) + (0x00 * 65536)
# 9  accent5  4472C4
$colors.Item(9).RGB  = 0x44 + (0x72 * 256) + (0xC4 * 65536)
# 10 accent6  70AD47
$colors.Item(10).RGB = 0x70 + (0xAD * 256) + (0x47 * 65536)
# 11 hlink    0563C1
$colors.Item(11).RGB = 0x05 + (0x63 * 256) + (0xC1 * 65536)
# 12 folHlink 954F72
$colors.Item(12).RGB = 0x95 + (0x4F * 256) + (0x72 * 65536)
